$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19 updates ---
$ws.Range("D19").Value = "one-dimensional temporal region"
$ws.Range("J19").Value = "intervention content and delivery"
$ws.Range("V19").Value = "PS"

# --- Row 21 updates ---
$ws.Range("D21").Value = "one-dimensional temporal region"
$ws.Range("J21").Value = "intervention content and delivery"

# --- Insert a new row at position 31 (pushes old rows 31..40 down to 32..41) ---
$ws.Rows(31).Insert()

# --- Give the whole new row the same "Proposed" white-fill look used by its
#     sibling rows (same visual treatment as rows whose Curation status is
#     "Proposed", e.g. the placebo intervention row right below it) ---
$ws.Range("A31:V31").Interior.Color = 16777215

# --- Populate the newly inserted row 31 with the BFO:0000038 entry ---
$ws.Range("A31").Value = "BFO:0000038"
$ws.Range("B31").Value = "one-dimensional temporal region"
$ws.Range("C31").Value = "A one-dimensional temporal region is a temporal region that is extended."
$ws.Range("D31").Value = "temporal region"
$ws.Range("J31").Value = "Intervention content and delivery"
$ws.Range("P31").Value = "LSR 1; LSR 2; LSR 3"
$ws.Range("Q31").Value = "Intervention content and delivery"
$ws.Range("S31").Value = "Proposed"
$ws.Range("V31").Value = "PS"
